# Update the "dSF" column (F) with repulled/pushed data and recalculated
# values. E ("dS0") is left untouched; F is updated to the newly pulled
# figures for each dated row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    3  = 1
    4  = -3
    5  = -2
    6  = 1
    7  = 3
    8  = 2
    9  = -1
    10 = -1
    11 = 1
    12 = -3
    13 = 1
    16 = -8
    17 = -6
    18 = -11
    19 = -10
    20 = -6
    21 = -1
    22 = 2
    23 = -2
    26 = -1
    28 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
